$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1
$ws.Range("I29").Value = 1
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 3
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 278

$ws.Range("H38").Value = 2963.3333
$ws.Range("I38").Value = 238.57143
$ws.Range("J38").Value = 12500
$ws.Range("K38").Value = 715.71429
$ws.Range("L38").Value = 37500
$ws.Range("M38").Value = -343.71429
$ws.Range("N38").Value = -38244

$ws.Range("H40").Value = 123820.16
$ws.Range("I40").Value = 3000000
$ws.Range("J40").Value = 3979.3333
$ws.Range("K40").Value = 3000000
$ws.Range("L40").Value = 3979.3333
$ws.Range("M40").Value = -2999825
$ws.Range("N40").Value = -4329.3333

$ws.Range("H58").Value = 2471.5
$ws.Range("I58").Value = 3600
$ws.Range("J58").Value = 1987.8572
$ws.Range("K58").Value = 10800
$ws.Range("L58").Value = 5963.571599999999
$ws.Range("M58").Value = -10650
$ws.Range("N58").Value = -6263.571599999999

$ws.Range("H69").Value = 9965.951999999999
$ws.Range("I69").Value = 9000
$ws.Range("J69").Value = 10014.25
$ws.Range("K69").Value = 27000
$ws.Range("L69").Value = 30042.75
$ws.Range("M69").Value = -26126
$ws.Range("N69").Value = -31790.75

$ws.Range("H72").Value = 9965.951999999999
$ws.Range("I72").Value = 9000
$ws.Range("J72").Value = 10014.25
$ws.Range("K72").Value = 81000
$ws.Range("L72").Value = 90128.25
$ws.Range("M72").Value = -76632
$ws.Range("N72").Value = -98864.25

$ws.Range("H86").Value = 100004520
$ws.Range("I86").Value = 90913384
$ws.Range("J86").Value = 111115910
$ws.Range("K86").Value = 90913384
$ws.Range("L86").Value = 111115910
$ws.Range("M86").Value = -90912261
$ws.Range("N86").Value = -111118156

$ws.Range("H87").Value = 170935
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 170935
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 170935
$ws.Range("N87").Value = -173431

$ws.Range("H89").Value = 100004520
$ws.Range("I89").Value = 90913384
$ws.Range("J89").Value = 111115910
$ws.Range("K89").Value = 454566920
$ws.Range("L89").Value = 555579550
$ws.Range("M89").Value = -454561304
$ws.Range("N89").Value = -555590782

$ws.Range("H90").Value = 170935
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 170935
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 512805
$ws.Range("N90").Value = -525285

$ws.Range("H96").Value = 2180.0527
$ws.Range("I96").Value = 1375.4445
$ws.Range("J96").Value = 2904.2
$ws.Range("K96").Value = 4126.333500000001
$ws.Range("L96").Value = 8712.599999999999
$ws.Range("M96").Value = -2753.333500000001
$ws.Range("N96").Value = -11458.6

$ws.Range("H112").Value = 3727
$ws.Range("I112").Value = 1993.5
$ws.Range("J112").Value = 3943.6875
$ws.Range("K112").Value = 5980.5
$ws.Range("L112").Value = 11831.0625
$ws.Range("M112").Value = -4872.5
$ws.Range("N112").Value = -14047.0625

$ws.Range("H132").Value = 2902.4807
$ws.Range("I132").Value = 2909.578
$ws.Range("J132").Value = 2856.8572
$ws.Range("K132").Value = 8728.734
$ws.Range("L132").Value = 8570.571599999999
$ws.Range("M132").Value = -6198.734
$ws.Range("N132").Value = -13630.5716

$ws.Range("H133").Value = 104398.6
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 104398.6
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 104398.6
$ws.Range("N133").Value = -114518.6

$ws.Range("H137").Value = 3390.1035
$ws.Range("I137").Value = 2591.25
$ws.Range("J137").Value = 5165.3335
$ws.Range("K137").Value = 7773.75
$ws.Range("L137").Value = 15496.0005
$ws.Range("M137").Value = -5223.75
$ws.Range("N137").Value = -20596.0005

$ws.Range("H138").Value = 2667.5762
$ws.Range("I138").Value = 1821
$ws.Range("J138").Value = 3334.5757
$ws.Range("K138").Value = 5463
$ws.Range("L138").Value = 10003.7271
$ws.Range("M138").Value = -323
$ws.Range("N138").Value = -20283.7271

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9693200
$ws.Range("I32").Value = 5556995
$ws.Range("J32").Value = 19238288
$ws.Range("K32").Value = 5556995
$ws.Range("L32").Value = 19238288
$ws.Range("M32").Value = -5556708
$ws.Range("N32").Value = -19238862

$ws.Range("H102").Value = 2088.111
$ws.Range("I102").Value = 1771
$ws.Range("J102").Value = 3198
$ws.Range("K102").Value = 1771
$ws.Range("L102").Value = 3198
$ws.Range("M102").Value = -149
$ws.Range("N102").Value = -6442

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 108999
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 108999
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 108999
$ws.Range("N108").Value = -116679

$ws.Range("H141").Value = 20354.5
$ws.Range("I141").Value = 20709

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3699.9033
$ws.Range("I31").Value = 1985.625
$ws.Range("J31").Value = 5528.467
$ws.Range("K31").Value = 1985.625
$ws.Range("L31").Value = 5528.467
$ws.Range("M31").Value = -1690.625
$ws.Range("N31").Value = -6118.467

$ws.Range("H34").Value = 3699.9033
$ws.Range("I34").Value = 1985.625
$ws.Range("J34").Value = 5528.467
$ws.Range("K34").Value = 1985.625
$ws.Range("L34").Value = 5528.467
$ws.Range("M34").Value = -1783.625
$ws.Range("N34").Value = -5932.467

$ws.Range("H132").Value = 4627.5
$ws.Range("I132").Value = 4404.2
$ws.Range("J132").Value = 4999.6665
$ws.Range("K132").Value = 13212.6
$ws.Range("L132").Value = 14998.9995
$ws.Range("M132").Value = -10682.6
$ws.Range("N132").Value = -20058.9995

$ws.Range("H134").Value = 4734.4116
$ws.Range("I134").Value = 3499.5
$ws.Range("J134").Value = 6498.5713
$ws.Range("K134").Value = 10498.5
$ws.Range("L134").Value = 19495.7139
$ws.Range("M134").Value = -7963.5
$ws.Range("N134").Value = -24565.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 89.583336
$ws.Range("I14").Value = 89.583336
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 268.750008
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -95.75000799999998

$ws.Range("H75").Value = 200
$ws.Range("I75").Value = 200
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 600
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 398
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 200
$ws.Range("I78").Value = 200
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 1800
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = 3192
$ws.Range("N78").ClearContents()

$ws.Range("H121").Value = 9093425
$ws.Range("I121").Value = 20000712
$ws.Range("J121").Value = 4020.1667
$ws.Range("K121").Value = 60002136
$ws.Range("L121").Value = 12060.5001
$ws.Range("M121").Value = -60000826
$ws.Range("N121").Value = -14680.5001

$ws.Range("H137").Value = 9368.308000000001
$ws.Range("I137").Value = 2241
$ws.Range("J137").Value = 15477.429
$ws.Range("K137").Value = 6723
$ws.Range("L137").Value = 46432.287
$ws.Range("M137").Value = -1623
$ws.Range("N137").Value = -56632.287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3037.25
$ws.Range("I80").Value = 1649
$ws.Range("J80").Value = 3500
$ws.Range("K80").Value = 1649
$ws.Range("L80").Value = 3500
$ws.Range("M80").Value = -651
$ws.Range("N80").Value = -5496

$ws.Range("H83").Value = 3037.25
$ws.Range("I83").Value = 1649
$ws.Range("J83").Value = 3500
$ws.Range("K83").Value = 8245
$ws.Range("L83").Value = 17500
$ws.Range("M83").Value = -3253
$ws.Range("N83").Value = -27484

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 14030.322
$ws.Range("I7").Value = 13122.542
$ws.Range("J7").Value = 17142.715
$ws.Range("K7").Value = 13122.542
$ws.Range("L7").Value = 17142.715
$ws.Range("M7").Value = -13010.542
$ws.Range("N7").Value = -17366.715

$ws.Range("H16").Value = 1521.9445
$ws.Range("I16").Value = 1226.3334
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 1226.3334
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -1056.3334
$ws.Range("N16").Value = -3340

$ws.Range("H46").Value = 3197.3333
$ws.Range("I46").Value = 1001
$ws.Range("J46").Value = 3636.6
$ws.Range("K46").Value = 1001
$ws.Range("L46").Value = 3636.6
$ws.Range("M46").Value = -813
$ws.Range("N46").Value = -4012.6

$ws.Range("H126").Value = 14030.322
$ws.Range("I126").Value = 13122.542
$ws.Range("J126").Value = 17142.715
$ws.Range("K126").Value = 39367.626
$ws.Range("L126").Value = 51428.145
$ws.Range("M126").Value = -36897.626
$ws.Range("N126").Value = -56368.145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2997.647
$ws.Range("I81").Value = 1975.25
$ws.Range("J81").Value = 3906.4443
$ws.Range("K81").Value = 3950.5
$ws.Range("L81").Value = 7812.8886
$ws.Range("M81").Value = -2889.5
$ws.Range("N81").Value = -9934.8886

$ws.Range("H84").Value = 2997.647
$ws.Range("I84").Value = 1975.25
$ws.Range("J84").Value = 3906.4443
$ws.Range("K84").Value = 19752.5
$ws.Range("L84").Value = 39064.443
$ws.Range("M84").Value = -14448.5
$ws.Range("N84").Value = -49672.443

$ws.Range("H107").Value = 623.5714
$ws.Range("I107").Value = 502
$ws.Range("J107").Value = 643.8333
$ws.Range("K107").Value = 1506
$ws.Range("L107").Value = 1931.4999
$ws.Range("M107").Value = 414
$ws.Range("N107").Value = -5771.4999

$ws.Range("H109").Value = 12000
$ws.Range("I109").Value = 12000
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 12000
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -10613
$ws.Range("N109").ClearContents()

$ws.Range("H126").Value = 4542.8887
$ws.Range("I126").Value = 4757.1763
$ws.Range("J126").Value = 900
$ws.Range("K126").Value = 14271.5289
$ws.Range("L126").Value = 2700
$ws.Range("M126").Value = -11801.5289
$ws.Range("N126").Value = -7640

$ws.Range("H136").Value = 2174.2144
$ws.Range("I136").Value = 1354.25
$ws.Range("J136").Value = 4224.125
$ws.Range("K136").Value = 4062.75
$ws.Range("L136").Value = 12672.375
$ws.Range("M136").Value = -1512.75
$ws.Range("N136").Value = -17772.375

Write-Host "done"
